# Auto-generated edit script: updates market-price derived cells (H:N)
# on multiple leve-profit worksheets, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1115.3846
$ws.Range("J40").Value = 1087.5
$ws.Range("L40").Value = 1087.5
$ws.Range("N40").Value = -1437.5
$ws.Range("H62").Value = 2972.9167
$ws.Range("I62").Value = 1739.2858
$ws.Range("J62").Value = 4700
$ws.Range("K62").Value = 1739.2858
$ws.Range("L62").Value = 4700
$ws.Range("M62").Value = -1115.2858
$ws.Range("N62").Value = -5948
$ws.Range("H65").Value = 2972.9167
$ws.Range("I65").Value = 1739.2858
$ws.Range("J65").Value = 4700
$ws.Range("K65").Value = 8696.429
$ws.Range("L65").Value = 23500
$ws.Range("M65").Value = -5576.429
$ws.Range("N65").Value = -29740
$ws.Range("H106").Value = 3636.182
$ws.Range("I106").Value = 2799.8
$ws.Range("K106").Value = 2799.8
$ws.Range("M106").Value = -2168.8
$ws.Range("H107").Value = 2347.04
$ws.Range("I107").Value = 1914.0952
$ws.Range("K107").Value = 1914.0952
$ws.Range("M107").Value = 5.904800000000023
$ws.Range("H116").Value = 424964.25
$ws.Range("I116").Value = 2003881.2
$ws.Range("J116").Value = 9459.789000000001
$ws.Range("K116").Value = 2003881.2
$ws.Range("L116").Value = 9459.789000000001
$ws.Range("M116").Value = -2000439.2
$ws.Range("N116").Value = -16343.789
$ws.Range("H137").Value = 2769.5469
$ws.Range("I137").Value = 1911.64
$ws.Range("J137").Value = 5833.5
$ws.Range("K137").Value = 5734.92
$ws.Range("L137").Value = 17500.5
$ws.Range("M137").Value = -3184.92
$ws.Range("N137").Value = -22600.5
$ws.Range("H138").Value = 4080.0435
$ws.Range("I138").Value = 662.05554
$ws.Range("J138").Value = 4911.446
$ws.Range("K138").Value = 1986.16662
$ws.Range("L138").Value = 14734.338
$ws.Range("M138").Value = 3153.83338
$ws.Range("N138").Value = -25014.338
$ws.Range("H141").Value = 3970.6
$ws.Range("I141").Value = 3955.257
$ws.Range("J141").Value = 4078
$ws.Range("K141").Value = 11865.771
$ws.Range("L141").Value = 12234
$ws.Range("M141").Value = -6685.771000000001
$ws.Range("N141").Value = -22594

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3894.5442
$ws.Range("I32").Value = 3613.5076
$ws.Range("J32").Value = 5463.6665
$ws.Range("K32").Value = 3613.5076
$ws.Range("L32").Value = 5463.6665
$ws.Range("M32").Value = -3326.5076
$ws.Range("N32").Value = -6037.6665
$ws.Range("H74").Value = 3861.8276
$ws.Range("I74").Value = 4125.6816
$ws.Range("K74").Value = 4125.6816
$ws.Range("M74").Value = -3251.6816
$ws.Range("H77").Value = 3861.8276
$ws.Range("I77").Value = 4125.6816
$ws.Range("K77").Value = 20628.408
$ws.Range("M77").Value = -16260.408
$ws.Range("H97").Value = 653.7241
$ws.Range("I97").Value = 714.0833
$ws.Range("J97").Value = 364
$ws.Range("K97").Value = 714.0833
$ws.Range("L97").Value = 364
$ws.Range("M97").Value = -218.0833
$ws.Range("N97").Value = -1356
$ws.Range("H132").Value = 2369.7827
$ws.Range("I132").Value = 1479.7941
$ws.Range("K132").Value = 4439.3823
$ws.Range("M132").Value = -1909.3823

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 16192.167
$ws.Range("I26").Value = 13430.6
$ws.Range("J26").Value = 30000
$ws.Range("K26").Value = 13430.6
$ws.Range("L26").Value = 30000
$ws.Range("M26").Value = -13138.6
$ws.Range("N26").Value = -30584
$ws.Range("H96").Value = 18990
$ws.Range("I96").Value = 7980
$ws.Range("J96").Value = 30000
$ws.Range("K96").Value = 7980
$ws.Range("L96").Value = 30000
$ws.Range("M96").Value = -5234
$ws.Range("N96").Value = -35492
$ws.Range("H134").Value = 1796.3827
$ws.Range("I134").Value = 1152.228
$ws.Range("J134").Value = 3326.25
$ws.Range("K134").Value = 3456.684
$ws.Range("L134").Value = 9978.75
$ws.Range("M134").Value = -921.6840000000002
$ws.Range("N134").Value = -15048.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8066931
$ws.Range("I31").Value = 1349.2
$ws.Range("J31").Value = 29417000
$ws.Range("K31").Value = 1349.2
$ws.Range("L31").Value = 29417000
$ws.Range("M31").Value = -1054.2
$ws.Range("N31").Value = -29417590
$ws.Range("H34").Value = 8066931
$ws.Range("I34").Value = 1349.2
$ws.Range("J34").Value = 29417000
$ws.Range("K34").Value = 1349.2
$ws.Range("L34").Value = 29417000
$ws.Range("M34").Value = -1147.2
$ws.Range("N34").Value = -29417404
$ws.Range("H99").Value = 7412180.5
$ws.Range("I99").Value = 11114633
$ws.Range("J99").Value = 7275.4443
$ws.Range("K99").Value = 11114633
$ws.Range("L99").Value = 7275.4443
$ws.Range("M99").Value = -11113135
$ws.Range("N99").Value = -10271.4443
$ws.Range("H126").Value = 7412180.5
$ws.Range("I126").Value = 11114633
$ws.Range("J126").Value = 7275.4443
$ws.Range("K126").Value = 33343899
$ws.Range("L126").Value = 21826.3329
$ws.Range("M126").Value = -33341429
$ws.Range("N126").Value = -26766.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1800.3334
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 2100.5
$ws.Range("K22").Value = 3600
$ws.Range("L22").Value = 6301.5
$ws.Range("M22").Value = -3431
$ws.Range("N22").Value = -6639.5
$ws.Range("H23").Value = 183.47058
$ws.Range("J23").Value = 249.4
$ws.Range("L23").Value = 748.2
$ws.Range("N23").Value = -1218.2
$ws.Range("H25").Value = 4420
$ws.Range("J25").Value = 4420
$ws.Range("L25").Value = 13260
$ws.Range("N25").Value = -13598
$ws.Range("H27").Value = 1800.3334
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 2100.5
$ws.Range("K27").Value = 3600
$ws.Range("L27").Value = 6301.5
$ws.Range("M27").Value = -3498
$ws.Range("N27").Value = -6505.5
$ws.Range("H30").Value = 4420
$ws.Range("J30").Value = 4420
$ws.Range("L30").Value = 13260
$ws.Range("N30").Value = -13464
$ws.Range("H113").Value = 595.34485
$ws.Range("I113").Value = 497.15555
$ws.Range("K113").Value = 1491.46665
$ws.Range("M113").Value = 678.5333499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10419132
$ws.Range("I80").Value = 14708266
$ws.Range("K80").Value = 14708266
$ws.Range("M80").Value = -14707268
$ws.Range("H83").Value = 10419132
$ws.Range("I83").Value = 14708266
$ws.Range("K83").Value = 73541330
$ws.Range("M83").Value = -73536338
$ws.Range("H97").Value = 993.2
$ws.Range("I97").Value = 934.3333
$ws.Range("J97").Value = 1081.5
$ws.Range("K97").Value = 934.3333
$ws.Range("L97").Value = 1081.5
$ws.Range("M97").Value = -438.3333
$ws.Range("N97").Value = -2073.5
$ws.Range("H132").Value = 2036.1356
$ws.Range("I132").Value = 1247.0244
$ws.Range("K132").Value = 3741.0732
$ws.Range("M132").Value = -1211.0732

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4169045.8
$ws.Range("I132").Value = 2193.4407
$ws.Range("J132").Value = 15875916
$ws.Range("K132").Value = 6580.3221
$ws.Range("L132").Value = 47627748
$ws.Range("M132").Value = -4050.3221
$ws.Range("N132").Value = -47632808
$ws.Range("H136").Value = 4351.8237
$ws.Range("I136").Value = 4611.8887
$ws.Range("K136").Value = 13835.6661
$ws.Range("M136").Value = -11285.6661

